$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.060.44"
$ws.Range("E2").Value = "  -0.45%  "

# Row 3
$ws.Range("D3").Value = "1.621.91"
$ws.Range("E3").Value = "  -1.49%  "

# Row 4
$ws.Range("E4").Value = "  +0.35%  "

# Row 5
$ws.Range("D5").Value = "'215.23"
$ws.Range("E5").Value = "  -0.82%  "

# Row 6
$ws.Range("E6").Value = "  +0.14%  "

# Row 7
$ws.Range("E7").Value = "  +0.25%  "

# Row 8
$ws.Range("E8").Value = "  -0.92%  "

# Row 9
$ws.Range("D9").Value = "'0.0624"

# Row 10
$ws.Range("D10").Value = "'20.13"

# Row 11
$ws.Range("E11").Value = "  -0.11%  "

# Row 12
$ws.Range("D12").Value = "1.630.13"
$ws.Range("E12").Value = "  -0.97%  "

# Row 13
$ws.Range("E13").Value = "  -0.72%  "

# Row 14
$ws.Range("E14").Value = "  -0.03%  "

# Row 15
$ws.Range("D15").Value = "27.048.98"
$ws.Range("E15").Value = "  -0.43%  "

# Row 16
$ws.Range("D16").Value = "'64.49"
$ws.Range("E16").Value = "  -4.63%  "

# Row 17
$ws.Range("E17").Value = "  -0.04%  "

# Row 18
$ws.Range("D18").Value = "'216.15"
$ws.Range("E18").Value = "  -1.34%  "

# Row 19
$ws.Range("D19").Value = "'1.01"
$ws.Range("E19").Value = "  +0.35%  "

# Row 20
$ws.Range("D20").Value = "'6.88"
$ws.Range("E20").Value = "  +0.85%  "

# Row 21
$ws.Range("E21").Value = "  -0.99%  "

# Row 23
$ws.Range("D23").Value = "'8.98"
$ws.Range("E23").Value = "  -2.51%  "

# Row 24
$ws.Range("D24").Value = "'147.66"
$ws.Range("E24").Value = "  -0.19%  "

# Row 25
$ws.Range("E25").Value = "  +0.28%  "

# Row 26
$ws.Range("D26").Value = "'7.27"
$ws.Range("E26").Value = "  -3.89%  "

# Row 27
$ws.Range("E27").Value = "  -0.42%  "

# Row 28
$ws.Range("D28").Value = "'15.57"
$ws.Range("E28").Value = "  -1.13%  "

# Row 29
$ws.Range("D29").Value = "'0.0503"
$ws.Range("E29").Value = "  -0.90%  "

# Row 30
$ws.Range("D30").Value = "'1.17"
$ws.Range("E30").Value = "  -0.80%  "

# Row 31
$ws.Range("E31").Value = "  -1.13%  "

# Row 32
$ws.Range("D32").Value = "'2.99"
$ws.Range("E32").Value = "  -1.45%  "

# Row 33
$ws.Range("D33").Value = "1.334.88"
$ws.Range("E33").Value = "  +5.55%  "

# Row 34
$ws.Range("E34").Value = "  -0.89%  "

# Row 35
$ws.Range("D35").Value = "'2.47"
$ws.Range("E35").Value = "  +0.37%  "

# Row 36
$ws.Range("D36").Value = "'0.0176"
$ws.Range("E36").Value = "  -1.01%  "

# Row 37
$ws.Range("D37").Value = "'0.542"
$ws.Range("E37").Value = "  -0.90%  "

# Row 38
$ws.Range("D38").Value = "'0.847"
$ws.Range("E38").Value = "  -0.15%  "

# Row 39
$ws.Range("E39").Value = "  +0.28%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.24"
$ws.Range("E40").Value = "  +0.49%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.801"
$ws.Range("E41").Value = "  -0.89%  "

# Row 42
$ws.Range("D42").Value = "'64.26"
$ws.Range("E42").Value = "  +3.67%  "

# Row 43
$ws.Range("D43").Value = "1.762.65"
$ws.Range("E43").Value = "  -1.36%  "

# Row 44
$ws.Range("D44").Value = "'5.21"
$ws.Range("E44").Value = "  -3.94%  "

# Row 45
$ws.Range("D45").Value = "'90.38"
$ws.Range("E45").Value = "  -1.17%  "

# Row 46
$ws.Range("D46").Value = "'1.60"
$ws.Range("E46").Value = "  -0.19%  "

# Row 47
$ws.Range("E47").Value = "  +22.36%  "

# Row 48
$ws.Range("E48").Value = "  -6.50%  "

# Row 49
$ws.Range("E49").Value = "  -0.09%  "

# Row 50
$ws.Range("D50").Value = "'0.0984"
$ws.Range("E50").Value = "  +1.12%  "

# Row 51
$ws.Range("D51").Value = "'7.54"
$ws.Range("E51").Value = "  -1.06%  "
